$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values may look numeric (single "." as decimal point); force
# text format so Excel keeps them as literal strings, matching the source data.
# Column E (Volume) values already carry padding spaces so they stay text naturally,
# but we set "@" there too for consistency/safety.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.647.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.15"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.22"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.598.63"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.628.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.69%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.43"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.68"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.90%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.210.25"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.64%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.502"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.29"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.769.07"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.03"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.67"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.61"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.89%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.20%  "
